$wb = $excel.ActiveWorkbook

# --- 1. "About" sheet: update the report date in C1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45387

# --- 2. "BAU Emissions" sheet: relabel row headers, update data, fix selection ---
$wsBau = $wb.Worksheets.Item("BAU Emissions")

# Relabel the "... : NoSettings" row headers to "... : test" (275 rows)
[void]$wsBau.Cells.Replace(" : NoSettings", " : test")

# Update the data series in row 94 (natural gas if / iron and steel) for 2032-2050
$wsBau.Range("M94").Value = 1001080
$wsBau.Range("N94").Value = 2002150
$wsBau.Range("O94").Value = 3003230
$wsBau.Range("P94").Value = 4004300
$wsBau.Range("Q94:AE94").Value = 5005380

# --- 3. Restore the selection on "BAU Emissions" to match the saved view state ---
$wsBau.Activate()
[void]$wsBau.Range("A30:AE280").Select()

# --- 4. Make "About" the active/selected tab, matching the saved view state ---
$wsAbout.Activate()
[void]$wsAbout.Range("E29").Select()
